# Apply updated cryptocurrency price / volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D2:D51 hold prices that look numeric (e.g. "1.028") but must stay
# plain text, exactly like the source data. Force text format first so the
# COM value-setter does not silently convert them into real numbers, then
# restore the default "Normal" style so the cells keep their original,
# unstyled appearance once the text is in place.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.596.34"
$ws.Range("D3").Value = "1.848.59"
$ws.Range("D4").Value = "1.028"
$ws.Range("D5").Value = "320.77"
$ws.Range("D6").Value = "1.025"
$ws.Range("D7").Value = "0.4371"
$ws.Range("D8").Value = "0.3743"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D10").Value = "0.8761"
$ws.Range("D11").Value = "21.48"
$ws.Range("D12").Value = "1.857.96"
$ws.Range("D13").Value = "5.495"
$ws.Range("D14").Value = "6.681"
$ws.Range("D15").Value = "0.07154"
$ws.Range("D16").Value = "82.82"
$ws.Range("D17").Value = "1.032"
$ws.Range("D18").Value = "0.000009016"
$ws.Range("D19").Value = "1.025"
$ws.Range("D21").Value = "27.619.82"
$ws.Range("D22").Value = "5.248"
$ws.Range("D24").Value = "2.070.82"
$ws.Range("D25").Value = "157.18"
$ws.Range("D26").Value = "1.931"
$ws.Range("D28").Value = "5.285"
$ws.Range("D29").Value = "1.954"
$ws.Range("D30").Value = "116.17"
$ws.Range("D31").Value = "0.09067"
$ws.Range("D33").Value = "0.7672"
$ws.Range("D35").Value = "2.874"
$ws.Range("D36").Value = "1.027"
$ws.Range("D38").Value = "0.01978"
$ws.Range("D39").Value = "0.05265"
$ws.Range("D40").Value = "0.5171"
$ws.Range("D41").Value = "2.799"
$ws.Range("D43").Value = "6.707"
$ws.Range("D44").Value = "8.576"
$ws.Range("D45").Value = "108.96"
$ws.Range("D46").Value = "10.55"
$ws.Range("D47").Value = "1.713"
$ws.Range("D48").Value = "0.4656"
$ws.Range("D49").Value = "0.06380"
$ws.Range("D50").Value = "1.892"

# Cells that already read as plain text (prices with two dots, like
# "27.596.34") are unaffected by the text-format cast above, so this is
# safe to apply uniformly across the whole price column.
$priceRange.Style = "Normal"

# Volume(1h) percentage strings (column E) are never numeric-looking, so
# they can be written directly without any text-format coercion.
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +2.69%  "
$ws.Range("E5").Value = "  +4.16%  "
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("E12").Value = "  -3.87%  "
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("E41").Value = "  +6.10%  "
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("E50").Value = "  +5.77%  "
$ws.Range("E51").Value = "  +6.23%  "
